# Apply the spreadsheet edit described by the commit:
# "Day4-first commit: The Website displays the list of Doctors matching
#  the searched condition"
#
# The sheet holds a single "number" column. The existing lookup value in
# A2 is replaced with a new phone/search number, and a new row (A3) is
# added containing an additional number, which also extends the used
# range of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Replace the existing search value in A2
$ws.Range("A2").Value = 9600794463

# Add the new number in the newly appended row A3
$ws.Range("A3").Value = 123456

# Leave A2 as the active/selected cell, matching the saved view state
$ws.Range("A2").Select()
